$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp footer
$ws.Range("A1").Value = "Datos actualizados a 15 de Agosto de 2020 a las 17:03"

# Update refreshed COVID-19 stats (and re-ranked country pairs) for the affected rows
# Row 4
$ws.Range("B4").Value = 5480231
$ws.Range("C4").Value = 3965
$ws.Range("D4").Value = 2876792
$ws.Range("E4").Value = 2431847
$ws.Range("G4").Value = 57
$ws.Range("H4").Value = 171592

# Row 6
$ws.Range("B6").Value = 2557342
$ws.Range("C6").Value = 32120
$ws.Range("D6").Value = 1835640
$ws.Range("E6").Value = 672221
$ws.Range("G6").Value = 347
$ws.Range("H6").Value = 49481

# Row 18
$ws.Range("D18").Value = 205697
$ws.Range("E18").Value = 71175
$ws.Range("G18").Value = 38
$ws.Range("H18").Value = 5565

# Row 22
$ws.Range("B22").Value = 223981
$ws.Range("C22").Value = 207
$ws.Range("E22").Value = 12142

# Row 41
$ws.Range("B41").Value = 75697
$ws.Range("C41").Value = 512
$ws.Range("D41").Value = 67519
$ws.Range("E41").Value = 7680
$ws.Range("G41").Value = 4
$ws.Range("H41").Value = 498

# Row 46
$ws.Range("B46").Value = 62313
$ws.Range("C46").Value = 885
$ws.Range("D46").Value = 50183
$ws.Range("E46").Value = 9775
$ws.Range("G46").Value = 14
$ws.Range("H46").Value = 2355

# Row 61
$ws.Range("A61").Value = "Uzbekistan"
$ws.Range("B61").Value = 34251
$ws.Range("C61").Value = 430
$ws.Range("D61").Value = 28661
$ws.Range("E61").Value = 5367
$ws.Range("G61").Value = 3
$ws.Range("H61").Value = 223

# Row 62
$ws.Range("A62").Value = "Azerbaiyan"
$ws.Range("B62").Value = 34107
$ws.Range("C62").Value = 89
$ws.Range("D62").Value = 31697
$ws.Range("E62").Value = 1904
$ws.Range("G62").Value = 2
$ws.Range("H62").Value = 506

# Row 64
$ws.Range("A64").Value = "Moldavia"
$ws.Range("B64").Value = 29905
$ws.Range("C64").Value = 422
$ws.Range("D64").Value = 20908
$ws.Range("E64").Value = 8102
$ws.Range("G64").Value = 11
$ws.Range("H64").Value = 895

# Row 65
$ws.Range("A65").Value = "Kenia"
$ws.Range("B65").Value = 29849
$ws.Range("C65").Value = 515
$ws.Range("D65").Value = 15970
$ws.Range("E65").Value = 13407
$ws.Range("G65").Value = 7
$ws.Range("H65").Value = 472

# Row 95
$ws.Range("B95").Value = 8029
$ws.Range("C95").Value = 40
$ws.Range("D95").Value = 6815
$ws.Range("E95").Value = 1150
$ws.Range("G95").Value = 1
$ws.Range("H95").Value = 64

# Row 142
$ws.Range("A142").Value = "Uganda"
$ws.Range("B142").Value = 1434
$ws.Range("C142").Value = 49
$ws.Range("D142").Value = 1142
$ws.Range("E142").Value = 279
$ws.Range("G142").Value = 1
$ws.Range("H142").Value = 13

# Row 143
$ws.Range("A143").Value = "Uruguay"
$ws.Range("B143").Value = 1421
$ws.Range("D143").Value = 1182
$ws.Range("E143").Value = 201
$ws.Range("H143").Value = 38

# Row 168
$ws.Range("A168").Value = "Trinidad yTobago"
$ws.Range("B168").Value = 474
$ws.Range("C168").Value = 48
$ws.Range("D168").Value = 139
$ws.Range("E168").Value = 325
$ws.Range("H168").Value = 10

# Row 169
$ws.Range("A169").Value = "Guadalupe"
$ws.Range("B169").Value = 446
$ws.Range("D169").Value = 289
$ws.Range("E169").Value = 143
$ws.Range("H169").Value = 14

# Row 213
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0

# Row 214
$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1
